$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the country names for rows 78 and 79 (Dinamarca's case count overtook
# El Salvador's, so the two countries trade places in the ranking).
$ws.Range("A78").Value = "Dinamarca"
$ws.Range("A79").Value = "El Salvador"

# Update the "last updated" timestamp string (cell A1).
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 15:04"

# Row 18 (Irak)
$ws.Range("B18").Value = 379141
$ws.Range("C18").Value = 3210
$ws.Range("D18").Value = 307482
$ws.Range("E18").Value = 62260
$ws.Range("G18").Value = 52
$ws.Range("H18").Value = 9399

# Row 20 (Arabia Saudita)
$ws.Range("B20").Value = 336387
$ws.Range("C20").Value = 390
$ws.Range("D20").Value = 321485
$ws.Range("E20").Value = 10027
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 4875

# Row 26 (Alemania)
$ws.Range("B26").Value = 300337
$ws.Range("C26").Value = 309
$ws.Range("D26").Value = 261900
$ws.Range("E26").Value = 28840

# Row 33 (Paises Bajos)
$ws.Range("B33").Value = 135892
$ws.Range("C33").Value = 4003
$ws.Range("G33").Value = 5
$ws.Range("H33").Value = 6454

# Row 36 (Catar)
$ws.Range("B36").Value = 126498
$ws.Range("C36").Value = 159
$ws.Range("D36").Value = 123475
$ws.Range("E36").Value = 2807

# Row 70 (Estado de Palestina)
$ws.Range("B70").Value = 41498
$ws.Range("C70").Value = 420
$ws.Range("D70").Value = 34698
$ws.Range("E70").Value = 6470
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 330

# Row 78 (now Dinamarca)
$ws.Range("B78").Value = 29680
$ws.Range("C78").Value = 378
$ws.Range("D78").Value = 22716
$ws.Range("E78").Value = 6306
$ws.Range("G78").Value = 4
$ws.Range("H78").Value = 658

# Row 79 (now El Salvador)
$ws.Range("B79").Value = 29358
$ws.Range("D79").Value = 24175
$ws.Range("E79").Value = 4320
$ws.Range("G79").Value = 6
$ws.Range("H79").Value = 863

# Row 102 (Consejo Danes para los Refugiados)
$ws.Range("B102").Value = 10760
$ws.Range("C102").Value = 8
$ws.Range("D102").Value = 10239
$ws.Range("E102").Value = 247

# Row 107 (Tayikistan)
$ws.Range("B107").Value = 9935
$ws.Range("C107").Value = 40
$ws.Range("D107").Value = 8749
$ws.Range("E107").Value = 1108
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 78

# Row 139 (Estonia)
$ws.Range("B139").Value = 3607
$ws.Range("C139").Value = 30
$ws.Range("E139").Value = 791
